# Fruta / hortaliza, semanal
# Insert a new weekly row at row 7 (existing rows 7-17 shift down to 8-18),
# then populate the new row 7 with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("7:7").Insert()

$ws.Range("A7").Value = 10
$ws.Range("B7").Value = "Vega Modelo de Temuco"
$ws.Range("C7").Value = "La Araucanía"
$ws.Range("D7").Value = 45044
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = 100112041
$ws.Range("G7").Value = "Fruto del paraíso"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 24000
$ws.Range("L7").Value = 24000
$ws.Range("M7").Value = 24000
$ws.Range("N7").Value = "`$/caja 18 kilos empedrada"
$ws.Range("O7").Value = "Región de Arica y Parinacota"
$ws.Range("P7").Value = 1333
$ws.Range("Q7").Value = 18
$ws.Range("R7").Value = "Hortaliza"
